# Update data: append 2020-04-17, 2020-04-18, 2020-04-19 rows (54-56) to the
# "Confirmados" and "Mortes" sheets.

$wb = $excel.ActiveWorkbook

$dates = @("2020-04-17", "2020-04-18", "2020-04-19")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

$confirmados = @{
    54 = @(135,110,370,1809,1059,2684,746,856,335,797,162,143,1021,557,195,874,2006,102,4349,463,802,92,164,926,12841,53,31);
    55 = @(142,132,393,1897,1193,3034,762,952,378,1040,171,161,1077,640,205,945,2193,123,4543,516,831,110,201,962,13894,71,33);
    56 = @(163,159,416,2044,1230,3252,827,1099,393,1205,174,168,1154,685,236,987,2459,145,4765,531,854,128,222,975,14267,83,33)
}

$mortes = @{
    54 = @(5,7,10,145,36,149,20,25,16,40,5,5,35,26,26,42,186,8,341,23,22,3,3,30,928,4,1);
    55 = @(5,7,10,161,37,176,24,28,18,44,5,5,39,33,26,46,205,9,387,24,24,3,3,31,991,5,1);
    56 = @(6,15,11,182,45,186,24,30,18,48,5,5,39,34,29,48,216,10,402,25,24,4,3,32,1015,5,1)
}

$sheetData = @{
    "Confirmados" = $confirmados;
    "Mortes" = $mortes
}

foreach ($sheetName in @("Confirmados", "Mortes")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $sheetData[$sheetName]

    for ($i = 0; $i -lt 3; $i++) {
        $r = 54 + $i
        $dateText = $dates[$i]

        # Write the date as literal text (shared string), not an auto-converted
        # date serial: enter it with a leading apostrophe (forces text) and then
        # reset the cell style back to Normal so no extra number format sticks.
        $cell = $ws.Range("A$r")
        $cell.Value = "'" + $dateText
        $cell.Style = "Normal"

        $vals = $rowsForSheet[$r]
        for ($j = 0; $j -lt $cols.Length; $j++) {
            $ws.Range($cols[$j] + "$r").Value = $vals[$j]
        }
    }
}

Write-Host "rows 54-56 added to Confirmados and Mortes"
